$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Add some unhappy case": swap the credential value in B2 for a new
# sample password, leaving the (now mismatched) mailto hyperlink as-is.
$ws.Range("B2").Value = "duy@6601"

# The other sample passwords that used to live in the shared-string table
# were only ever referenced by B2, so they fall out of the workbook once
# B2 no longer points at them - nothing further to do for those.

# The unused built-in "Hyperlink" cell style (it was never actually
# applied to B2) gets cleaned up too.
$wb.Styles.Item("Hyperlink").Delete()

# Mirror the author's trailing selection below the data.
$ws.Range("A3:XFD14").Select()
